$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 333 (before XK), shifting XK..ZW down by one.
$ws.Rows(333).Insert()

# Populate the new row with the WS (Samoa) mapping -> OCE region.
$ws.Range("A333").Value = "WS"
$ws.Range("B333").Value = "WS"
$ws.Range("C333").Value = "OCE"
$ws.Range("D333").Value = "OCE"
$ws.Range("E333").Value = "OCE"

# Restore the view state to match what was recorded after the edit
# (scrolled down near the newly inserted/edited rows).
$ws.Range("B319").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 311
